# Daily update at 8 AM UTC
# Appends the next day's data row to the "Wins Over Time" tracking sheet.
# The previously-last row (row 72, the "today" row) loses its highlighted
# date format and becomes a normal historical row, while the newly added
# row 73 becomes the new "today" row with the highlighted date format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 72 was the "current day" row (date-only format). Now that a new day
# has arrived, revert it to the standard historical-row format.
$ws.Range("A72").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Add the new day's data in row 73.
$ws.Range("A73").Value = 45660
$ws.Range("B73").Value = 169
$ws.Range("C73").Value = 167
$ws.Range("D73").Value = 171

# Mark row 73 as the new "current day" row with the date-only format.
$ws.Range("A73").NumberFormat = "YYYY-MM-DD"
